$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns E:G (meandam / lowerdam / upperdam), matching D1's style ---
$ws.Range("D1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)
$ws.Range("E1").Value = "meandam"
$ws.Range("F1").Value = "lowerdam"
$ws.Range("G1").Value = "upperdam"

# --- Per-row data: B (rp bucket code), C (curve -> now same as asset_type),
#     and new E/F/G mean/lower/upper damage values ---
$ws.Range("B2").Value = "rp0001"
$ws.Range("C2").Value = "plant"
$ws.Range("E2").Value = 544840920.5143894
$ws.Range("F2").Value = 408630690.385792
$ws.Range("G2").Value = 681051150.6429868
$ws.Range("B3").Value = "rp0001"
$ws.Range("C3").Value = "substation"
$ws.Range("E3").Value = 5989529.050910485
$ws.Range("F3").Value = 4492146.788182863
$ws.Range("G3").Value = 7486911.313638107
$ws.Range("B4").Value = "rp0002"
$ws.Range("C4").Value = "plant"
$ws.Range("E4").Value = 568529545.9378448
$ws.Range("F4").Value = 426397159.4533835
$ws.Range("G4").Value = 710661932.422306
$ws.Range("B5").Value = "rp0002"
$ws.Range("C5").Value = "substation"
$ws.Range("E5").Value = 6104699.094479713
$ws.Range("F5").Value = 4578524.320859784
$ws.Range("G5").Value = 7630873.868099641
$ws.Range("B6").Value = "rp0005"
$ws.Range("C6").Value = "plant"
$ws.Range("E6").Value = 784498565.2783953
$ws.Range("F6").Value = 588373923.9587964
$ws.Range("G6").Value = 980623206.5979941
$ws.Range("B7").Value = "rp0005"
$ws.Range("C7").Value = "substation"
$ws.Range("E7").Value = 8488540.979717484
$ws.Range("F7").Value = 6366405.734788112
$ws.Range("G7").Value = 10610676.22464686
$ws.Range("B8").Value = "rp0010"
$ws.Range("C8").Value = "plant"
$ws.Range("E8").Value = 838593735.9322069
$ws.Range("F8").Value = 628945301.9491551
$ws.Range("G8").Value = 1048242169.915259
$ws.Range("B9").Value = "rp0010"
$ws.Range("C9").Value = "substation"
$ws.Range("E9").Value = 8961232.731651062
$ws.Range("F9").Value = 6720924.548738295
$ws.Range("G9").Value = 11201540.91456383
$ws.Range("B10").Value = "rp0025"
$ws.Range("C10").Value = "plant"
$ws.Range("E10").Value = 909057474.6268866
$ws.Range("F10").Value = 681793105.9701649
$ws.Range("G10").Value = 1136321843.283608
$ws.Range("B11").Value = "rp0025"
$ws.Range("C11").Value = "substation"
$ws.Range("E11").Value = 9254578.070568452
$ws.Range("F11").Value = 6940933.552926338
$ws.Range("G11").Value = 11568222.58821057
$ws.Range("B12").Value = "rp0050"
$ws.Range("C12").Value = "plant"
$ws.Range("E12").Value = 959560619.0910087
$ws.Range("F12").Value = 719670464.3182564
$ws.Range("G12").Value = 1199450773.863761
$ws.Range("B13").Value = "rp0050"
$ws.Range("C13").Value = "substation"
$ws.Range("E13").Value = 9343986.736297682
$ws.Range("F13").Value = 7007990.052223261
$ws.Range("G13").Value = 11679983.4203721
$ws.Range("B14").Value = "rp0100"
$ws.Range("C14").Value = "plant"
$ws.Range("E14").Value = 1011973472.918356
$ws.Range("F14").Value = 758980104.6887666
$ws.Range("G14").Value = 1264966841.147944
$ws.Range("B15").Value = "rp0100"
$ws.Range("C15").Value = "substation"
$ws.Range("E15").Value = 9523796.86459252
$ws.Range("F15").Value = 7142847.64844439
$ws.Range("G15").Value = 11904746.08074065
$ws.Range("B16").Value = "rp0250"
$ws.Range("C16").Value = "plant"
$ws.Range("E16").Value = 1102696240.991511
$ws.Range("F16").Value = 827022180.7436334
$ws.Range("G16").Value = 1378370301.239389
$ws.Range("B17").Value = "rp0250"
$ws.Range("C17").Value = "substation"
$ws.Range("E17").Value = 9845210.531062711
$ws.Range("F17").Value = 7383907.898297034
$ws.Range("G17").Value = 12306513.16382839
$ws.Range("B18").Value = "rp0500"
$ws.Range("C18").Value = "plant"
$ws.Range("E18").Value = 1159107516.433252
$ws.Range("F18").Value = 869330637.3249388
$ws.Range("G18").Value = 1448884395.541565
$ws.Range("B19").Value = "rp0500"
$ws.Range("C19").Value = "substation"
$ws.Range("E19").Value = 10117982.16083502
$ws.Range("F19").Value = 7588486.620626263
$ws.Range("G19").Value = 12647477.70104377
$ws.Range("B20").Value = "rp1000"
$ws.Range("C20").Value = "plant"
$ws.Range("E20").Value = 1214875594.003349
$ws.Range("F20").Value = 911156695.5025115
$ws.Range("G20").Value = 1518594492.504186
$ws.Range("B21").Value = "rp1000"
$ws.Range("C21").Value = "substation"
$ws.Range("E21").Value = 10494878.79046413
$ws.Range("F21").Value = 7871159.092848093
$ws.Range("G21").Value = 13118598.48808016
